# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" table (rows 16-23, columns B:J) is re-sorted so the
# periods run in ascending chronological order (2010, 2011, 2012, 2101,
# 2102, 2103, 2104, 2105) instead of the previous descending order. Only
# columns E (Periodo Mora) and F (Valor Mora) actually differ row to row,
# so update those two columns in place to match the new ascending order.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$periodos = @("2010", "2011", "2012", "2101", "2102", "2103", "2104", "2105")
$valores  = @(35112, 35112, 35112, 35112, 35112, 35112, 35112, 28090)

for ($i = 0; $i -lt $periodos.Length; $i++) {
    $row = 16 + $i
    $ws.Cells.Item($row, 5).Value = $periodos[$i]
    $ws.Cells.Item($row, 6).Value = $valores[$i]
}
